# Auto-update TW Market Data: 2026-02-17 07:52:14 UTC
#
# The underlying daily-signals feed dropped "2308 / 台達電 / Delta Elec"
# (its %_Vol_vs_TB ranking fell out of the tracked list). Remove it from
# both the Daily Signals and 21-Day Trend sheets (rows shift up), then
# refresh the Industry Analysis aggregate for the "Power Supply" sector
# (now only Lite-On remains) and the recalculated MACD values on the
# Favorites sheet.

$wb = $excel.ActiveWorkbook

# --- 1_Daily_Signals: remove the Delta Elec (2308) row ---------------------
$ws1 = $wb.Worksheets.Item("1_Daily_Signals")
$ws1.Rows.Item(17).Delete()

# --- 2_21Day_Trend: remove the Delta Elec (2308) row ------------------------
$ws2 = $wb.Worksheets.Item("2_21Day_Trend")
$ws2.Rows.Item(7).Delete()

# --- 3_Industry_Analysis: Power Supply sector now has only Lite-On ---------
# (Delta Elec dropped out, so the sector average/sum/stock-count shrink and
# the row sinks from rank 11 down to rank 14; rows in between shift up.)
$ws3 = $wb.Worksheets.Item("3_Industry_Analysis")

$ws3.Range("A11").Value = "Design Service"
$ws3.Range("B11").Value = 9.855
$ws3.Range("C11").Value = 0.61
$ws3.Range("D11").Value = 1.624
$ws3.Range("E11").Value = 2

$ws3.Range("A12").Value = "Compound Semi"
$ws3.Range("B12").Value = 8.293333333333333
$ws3.Range("C12").Value = 0.9966666666666667
$ws3.Range("D12").Value = 10.299
$ws3.Range("E12").Value = 3

$ws3.Range("A13").Value = "Steel"
$ws3.Range("B13").Value = 7.73
$ws3.Range("C13").Value = 0.54
$ws3.Range("D13").Value = 1.671
$ws3.Range("E13").Value = 1

$ws3.Range("A14").Value = "Power Supply"
$ws3.Range("B14").Value = 7.49
$ws3.Range("C14").Value = 1.05
$ws3.Range("D14").Value = 4.867
$ws3.Range("E14").Value = 1

# --- 4_My_Favorites: refreshed MACD (column H) values -----------------------
$ws4 = $wb.Worksheets.Item("4_My_Favorites")

$ws4.Range("H2").Value = 80.2774
$ws4.Range("H3").Value = 94.9751
$ws4.Range("H4").Value = 3.7306
$ws4.Range("H5").Value = -2.1749
$ws4.Range("H6").Value = 5.9959
$ws4.Range("H7").Value = -1.7388
